$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new row of daily stats (Daily update at 8 AM UTC)
$ws.Range("A43").Value = 45992
$ws.Range("B43").Value = 97
$ws.Range("C43").Value = 108
$ws.Range("D43").Value = 103

# Match the date-style formatting used in column A for the rest of the data
$ws.Range("A43").NumberFormat = $ws.Range("A42").NumberFormat
